$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table per row: Fecha(D), Calidad(L), Volumen(M), Precio minimo(N),
# Precio maximo(O), Precio promedio ponderado(P), Unidad(Q), Precio $/Kg(S), Kg/unidad(T)
$rowsData = @{
    385 = @(45041, "Especial", 200, 12000, 12000, 12000, "$/bandeja 18 kilos", 667, 18)
    386 = @(45041, "Primera", 250, 10000, 10000, 10000, "$/bandeja 18 kilos", 556, 18)
    387 = @(44705, "Especial", 100, 10000, 10000, 10000, "$/bandeja 18 kilos", 556, 18)
    388 = @(44705, "Primera", 150, 8000, 8000, 8000, "$/bandeja 18 kilos", 444, 18)
    389 = @(44705, "Segunda", 130, 6000, 6000, 6000, "$/bandeja 18 kilos", 333, 18)
    390 = @(44448, "Especial", 180, 13000, 13000, 13000, "$/caja 18 kilos", 722, 18)
    391 = @(44448, "Primera", 230, 11000, 11000, 11000, "$/caja 18 kilos", 611, 18)
    392 = @(44342, "Primera", 250, 8000, 9000, 8600, "$/bandeja 18 kilos", 478, 18)
    393 = @(44727, "Especial", 100, 10000, 10000, 10000, "$/bandeja 18 kilos", 556, 18)
    394 = @(44727, "Primera", 250, 8000, 8000, 8000, "$/bandeja 18 kilos", 444, 18)
    395 = @(44391, "Especial", 200, 11000, 11000, 11000, "$/bandeja 18 kilos", 611, 18)
    396 = @(44391, "Primera", 250, 9500, 9500, 9500, "$/bandeja 18 kilos", 528, 18)
    397 = @(44875, "Especial", 250, 12000, 12000, 12000, "$/caja 10 kilos", 1200, 10)
    398 = @(44875, "Primera", 200, 10000, 10000, 10000, "$/caja 10 kilos", 1000, 10)
    399 = @(44454, "Primera", 100, 12000, 12000, 12000, "$/bandeja 18 kilos", 667, 18)
    400 = @(44454, "Segunda", 80, 10000, 10000, 10000, "$/bandeja 18 kilos", 556, 18)
    401 = @(44426, "Especial", 100, 12000, 12000, 12000, "$/caja 18 kilos", 667, 18)
    402 = @(44426, "Primera", 120, 10000, 10000, 10000, "$/caja 18 kilos", 556, 18)
}

foreach ($r in 385..400) {
    $vals = $rowsData[$r]
    $ws.Range("D$r").Value2 = $vals[0]
    $ws.Range("L$r").Value2 = $vals[1]
    $ws.Range("M$r").Value2 = $vals[2]
    $ws.Range("N$r").Value2 = $vals[3]
    $ws.Range("O$r").Value2 = $vals[4]
    $ws.Range("P$r").Value2 = $vals[5]
    $ws.Range("Q$r").Value2 = $vals[6]
    $ws.Range("S$r").Value2 = $vals[7]
    $ws.Range("T$r").Value2 = $vals[8]
}

# Append two new rows (401, 402): fixed columns A,B,C,E,F,G,H,I,J,K,R copied
# from row 400 (identical across the whole block); D,L,M,N,O,P,Q,S,T from table.
foreach ($r in 401..402) {
    $ws.Range("A$r").Value2 = $ws.Range("A400").Value2
    $ws.Range("B$r").Value2 = $ws.Range("B400").Value2
    $ws.Range("C$r").Value2 = $ws.Range("C400").Value2
    $vals = $rowsData[$r]
    $ws.Range("D$r").Value2 = $vals[0]
    $ws.Range("D$r").NumberFormat = $ws.Range("D400").NumberFormat
    $ws.Range("E$r").Value2 = $ws.Range("E400").Value2
    $ws.Range("F$r").Value2 = $ws.Range("F400").Value2
    $ws.Range("G$r").Value2 = $ws.Range("G400").Value2
    $ws.Range("H$r").Value2 = $ws.Range("H400").Value2
    $ws.Range("I$r").Value2 = $ws.Range("I400").Value2
    $ws.Range("J$r").Value2 = $ws.Range("J400").Value2
    $ws.Range("K$r").Value2 = $ws.Range("K400").Value2
    $ws.Range("L$r").Value2 = $vals[1]
    $ws.Range("M$r").Value2 = $vals[2]
    $ws.Range("N$r").Value2 = $vals[3]
    $ws.Range("O$r").Value2 = $vals[4]
    $ws.Range("P$r").Value2 = $vals[5]
    $ws.Range("Q$r").Value2 = $vals[6]
    $ws.Range("R$r").Value2 = $ws.Range("R400").Value2
    $ws.Range("S$r").Value2 = $vals[7]
    $ws.Range("T$r").Value2 = $vals[8]
}